$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1434.7237589384115
$ws.Range("E2").Value = 1727.5359942936627
$ws.Range("J2").Value = 471.96269851437467

$ws.Range("C3").Value = 1090.0368042981816
$ws.Range("E3").Value = 1761.269001654387
$ws.Range("J3").Value = 534.561144794189

$ws.Range("C4").Value = 1091.3452497379058
$ws.Range("E4").Value = 1859.0263156237197
$ws.Range("J4").Value = 849.0588754685059

$ws.Range("C5").Value = 1084.9869206170813
$ws.Range("E5").Value = 1883.4933627092037
$ws.Range("J5").Value = 719.2840458466168

$ws.Range("C6").Value = 583.9977800533463
$ws.Range("E6").Value = 1985.1420231656718
$ws.Range("J6").Value = 356.53577602764267

$ws.Range("C7").Value = 818.0998139706464
$ws.Range("E7").Value = 1729.864018678907
$ws.Range("J7").Value = 532.4200212299824

$ws.Range("C8").Value = 939.2589368041106
$ws.Range("E8").Value = 1705.2557394324979
$ws.Range("J8").Value = 647.3504003121201

$ws.Range("C9").Value = 993.5273242633199
$ws.Range("E9").Value = 1762.7113861234513
$ws.Range("J9").Value = 951.9623051755061

$ws.Range("C10").Value = 1057.7232374961402
$ws.Range("E10").Value = 2045.9363133819963
$ws.Range("J10").Value = 727.0181732426795

$ws.Range("C11").Value = 1208.5769622812456
$ws.Range("E11").Value = 1733.641494546538
$ws.Range("J11").Value = 533.938009398067
